# US-3925: New TCs for contract-duration legal clause + 24/36 month configs
# Updates D02Variables, ProfessionalInternet and EnterpriseVoice sheets.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) ProfessionalInternet (row 4 gets a value, row 5 is brand new)
#    Written first so the "36 Months VDSL2" string lands at sst#235,
#    matching the order in which the strings were first introduced.
# ------------------------------------------------------------------
$wsPI = $wb.Worksheets.Item("ProfessionalInternet")

$wsPI.Range("B4").ClearFormats()
$wsPI.Range("B4").Value = "New,36 Months,VDSL2,--None--,NotApplicable,NotApplicable,Internet Pro,NotApplicable"

# ------------------------------------------------------------------
# 2) D02Variables: new legal-clause variables (EN/FR/NL)
# ------------------------------------------------------------------
$wsD02 = $wb.Worksheets.Item("D02Variables")

$wsD02.Range("B2").Value = "By signing this order form, the Customer with fewer than 50 employees who opts for a contract of more than 24 months acknowledges having been duly informed of his right to subscribe to a contract of a shorter duration and of the related commercial conditions and declares that he renounces it."
$wsD02.Range("B3").Value = "En signant le présent bon de commande, le Client de moins de 50 employés qui opte pour un contrat d’une durée supérieure à 24 mois reconnaît avoir été dûment informé de son droit de souscrire à un contrat d’une durée inférieure et des conditions commerciales y relatives et déclare y renoncer."
$wsD02.Range("B4").Value = "Door ondertekening van deze bestelbon erkent de Klant met minder dan 50 werknemers die opteert voor een contract van meer dan 24 maanden, naar behoren te zijn ingelicht over zijn recht om in te tekenen op een contract van kortere duur en de daaraan verbonden commerciële voorwaarden en verklaart hij daarvan af te zien."

$wsD02.Range("A2").Value = "telcoLegalClauseAgreeDoc_EN"
$wsD02.Range("A3").Value = "telcoLegalClauseAgreeDoc_FR"
$wsD02.Range("A4").Value = "telcoLegalClauseAgreeDoc_NL"

# ------------------------------------------------------------------
# 3) EnterpriseVoice (row 3 gets a value)
# ------------------------------------------------------------------
$wsEV = $wb.Worksheets.Item("EnterpriseVoice")

$wsEV.Range("B3").ClearFormats()
$wsEV.Range("B3").HorizontalAlignment = -4131
$wsEV.Range("B3").Value = "New,36 Months,GPON,Voice Only,Existing,Bizz IP Box Other IPBX,120,VoIP Individual number,5"

# ------------------------------------------------------------------
# 4) ProfessionalInternet row 4 (col A) + second new row (24-months config)
#    "contractHigher24MonthsConfiguration" is introduced here (reused
#    as-is on EnterpriseVoice!A3 right after).
# ------------------------------------------------------------------
$wsPI.Range("A4").Value = "contractHigher24MonthsConfiguration"
$wsEV.Range("A3").Value = "contractHigher24MonthsConfiguration"

$wsPI.Range("A5").Value = "contractEqual24MonthsConfiguration"
$wsPI.Range("B5").Value = "New,24 Months,VDSL2,--None--,NotApplicable,NotApplicable,Internet Pro,NotApplicable"

# ------------------------------------------------------------------
# Column widths (best fit) -- approximate the author's AutoFit result
# ------------------------------------------------------------------
$wsPI.Columns.Item(1).ColumnWidth = 43.2857142857143
$wsEV.Columns.Item(1).ColumnWidth = 43.2857142857143

# ------------------------------------------------------------------
# Selections / active cells per sheet
# ------------------------------------------------------------------
$wsD02.Range("A24").Select()

$wsPI.Range("B11").Select()

$wsD03 = $wb.Worksheets.Item("D03NonQuotableProducts")
$wsD03.Range("B19").Select()

# EnterpriseVoice becomes the active / selected sheet, cell A2 selected.
$wsEV.Activate()
$wsEV.Range("A2").Select()
